$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Add MAX:/AVERAGE/MIN: labels in the order that reproduces the
# --- original author's shared-string table ordering (4=MAX:, 5=AVERAGE, 6=MIN:)
$ws.Range("A21").Value() = "MAX:"
$ws.Range("J3").Value() = "AVERAGE"
$ws.Range("A22").Value() = "MIN:"

# --- New column J: per-row AVERAGE(B:I) for rows 4..20
for ($r = 4; $r -le 20; $r++) {
    $ws.Cells.Item($r, 10).Formula() = "=AVERAGE(B$r`:I$r)"
}

# --- New row 21 (MAX) and row 22 (MIN) across columns B..J
for ($c = 2; $c -le 10; $c++) {
    $col = [char](64 + $c)
    $ws.Cells.Item(21, $c).Formula() = "=MAX($col`4:$col`20)"
    $ws.Cells.Item(22, $c).Formula() = "=MIN($col`4:$col`20)"
}

# --- Update the saved selection to match the target workbook
[void]$ws.Range("A17").Select()
